# Updated main GSC export data.
#
# The export window rolls forward by one day: the oldest day (2025-11-01)
# drops off the top of the "Chart" table and three new days
# (2026-01-27, 2026-01-28, 2026-01-29) are appended at the bottom. Because
# the window shifted, nearly every remaining row keeps the same
# "No video indexed" / "Video indexed" / "Impressions" values it already
# had for its calendar date - only the very last couple of days carry the
# freshly-exported numbers. The "Table" sheet's failure count is refreshed
# to match the new last-row total.

$wb = $excel.ActiveWorkbook
$chart = $wb.Worksheets.Item("Chart")
$table = $wb.Worksheets.Item("Table")

# Helper: force a plain literal-text cell value (shared string), bypassing
# Excel's automatic "this looks like a date/number" inference that a plain
# `.Value = "2026-01-27"` assignment would trigger. We park the text behind
# a quoted formula, then paste the computed result back over itself as a
# value only - this also strips the formula while leaving the cell's style
# untouched (no NumberFormat fiddling, so no new style entries are added).
function Set-TextValue {
    param($cell, [string]$text)
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

# Drop the oldest day (row 2, 2025-11-01) - this shifts every following row
# up by one, so row 3 (2025-11-02) becomes the new row 2, etc. All the
# daily counts travel with their original date, unchanged.
$chart.Rows.Item(2).Delete()

# After the delete, the last populated row is 87 (2026-01-26). Append the
# three newly exported days.
Set-TextValue $chart.Cells.Item(88, 1) "2026-01-27"
$chart.Cells.Item(88, 2).Value = 20
$chart.Cells.Item(88, 3).Value = 1
$chart.Cells.Item(88, 4).Value = 1

Set-TextValue $chart.Cells.Item(89, 1) "2026-01-28"
$chart.Cells.Item(89, 2).Value = 20
$chart.Cells.Item(89, 3).Value = 1
$chart.Cells.Item(89, 4).Value = 0

Set-TextValue $chart.Cells.Item(90, 1) "2026-01-29"
$chart.Cells.Item(90, 2).Value = 20
$chart.Cells.Item(90, 3).Value = 1
# Impressions for the freshly-added final day aren't reported yet - leave
# it blank (matches the export's empty placeholder for in-progress days).
$chart.Cells.Item(90, 4).Value = ""

# The validation/failure summary on the "Table" sheet tracks the latest
# "No video indexed" count, which is now 20 (was 19).
$table.Cells.Item(2, 3).Value = 20
